$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the anova/df/F/P helper labels in column E (rows 6-9) ---
# (their shared-string entries disappear from sharedStrings.xml once unreferenced)
$ws.Range("E6:E9").ClearContents()

# --- Mark B2:C6 with an explicit "No Fill" so the cells carry an applied style ---
$ws.Range("B2:C6").Interior.ColorIndex = -4142

# --- Update the second table's data (rows 17-27, columns B & C) ---
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 6
$ws.Range("B18").Value = 10
$ws.Range("C18").Value = 10
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = 9
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 8
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 5
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 4
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 8
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 6
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 5
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = 6
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 7

# --- Update the sheet's current selection to match the saved view ---
[void]$ws.Range("A26:XFD26").Select()
